$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 326, pushing existing rows 326:361 down to 327:362
$ws.Rows(326).Insert()

# Populate the new row 326 with the new price-report entry
$ws.Cells.Item(326, 1).Value = 10
$ws.Cells.Item(326, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(326, 3).Value = "La Araucanía"
$ws.Cells.Item(326, 4).Value = 44858
$ws.Cells.Item(326, 5).Value = 9
$ws.Cells.Item(326, 6).Value = 100114013
$ws.Cells.Item(326, 7).Value = "Zanahoria"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 530
$ws.Cells.Item(326, 11).Value = 20000
$ws.Cells.Item(326, 12).Value = 22000
$ws.Cells.Item(326, 13).Value = 20943
$ws.Cells.Item(326, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(326, 15).Value = "Región del Maule"
$ws.Cells.Item(326, 16).Value = 1047
$ws.Cells.Item(326, 17).Value = 20
$ws.Cells.Item(326, 18).Value = "Hortaliza"
